$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C10) from serial date 45183 to 45184
$ws.Range("C2:C10").Value = 45184
